# Gear Ratio.xlsx - add "Hemera" and "Mellow NF Cannon" rows, switch the
# circumference formula from a hard-coded Pi approximation to PI().

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ComparoGearing")

# --- 1. Swap the F column's hard-coded 3.1416 constant for PI() -----------
# (F2:F5 keeps its shared-formula shape; just change the multiplier.)
$ws.Range("F2:F5").Formula = "=E2*PI()"

# --- 2. Row 6: Hemera -------------------------------------------------------
$ws.Range("A6").Value = "Hemera"
$ws.Range("D6").Value = 3.32
$ws.Range("E6").Value = 8.27
$ws.Range("F6").Formula = "=E6*PI()"
$ws.Range("G6").Formula = "=1/F6"
$ws.Range("H6").Formula = "=G6*D6"
$ws.Range("I6").Formula = "=F6/D6"
$ws.Range("J6").Formula = "=H6*200*16"
$ws.Range("K6").Formula = "=1/H6*(1.75/2)*(1.75/2)*3.1416*8"

# --- 3. Row 7: Mellow NF Cannon --------------------------------------------
$ws.Range("A7").Value = "Mellow NF Cannon"
$ws.Range("B7").Value = 39
$ws.Range("C7").Value = 2
$ws.Range("D7").Formula = "=B7/C7"
$ws.Range("E7").Value = 18.8
$ws.Range("F7").Formula = "=E7*PI()"
$ws.Range("G7").Formula = "=1/F7"
$ws.Range("H7").Formula = "=G7*D7"
$ws.Range("I7").Formula = "=F7/D7"
$ws.Range("J7").Formula = "=H7*200*16"
$ws.Range("K7").Formula = "=1/H7*(1.75/2)*(1.75/2)*3.1416*8"

# --- 4. Formatting: extend the bordered-table look down into the new rows --
$ws.Range("A5:K5").Copy()
$ws.Range("A6:K6").PasteSpecial(-4122)
$ws.Range("A5:K5").Copy()
$ws.Range("A7:K7").PasteSpecial(-4122)
$ws.Range("J7:K7").NumberFormat = $ws.Range("J6").NumberFormat()
$excel.CutCopyMode = 0

Write-Output "done"
